$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310.4516
$ws.Range("I33").Value = 319.13333
$ws.Range("K33").Value = 319.13333
$ws.Range("M33").Value = -90.13333
$ws.Range("H107").Value = 21300
$ws.Range("I107").Value = 100
$ws.Range("J107").Value = 42500
$ws.Range("K107").Value = 100
$ws.Range("L107").Value = 42500
$ws.Range("M107").Value = 1820
$ws.Range("N107").Value = -46340
$ws.Range("H116").Value = 1987.5
$ws.Range("I116").Value = 1987.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1987.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1454.5
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value = 1476.4584
$ws.Range("I137").Value = 1113.0588
$ws.Range("J137").Value = 2359
$ws.Range("K137").Value = 3339.1764
$ws.Range("L137").Value = 7077
$ws.Range("M137").Value = -789.1764000000003
$ws.Range("N137").Value = -12177
$ws.Range("H138").Value = 1511.4556
$ws.Range("I138").Value = 1019.72546
$ws.Range("J138").Value = 2154.487
$ws.Range("K138").Value = 3059.17638
$ws.Range("L138").Value = 6463.461
$ws.Range("M138").Value = 2080.82362
$ws.Range("N138").Value = -16743.461
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2833.3333
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 3500
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = -1387
$ws.Range("N2").Value = -3726
$ws.Range("H32").Value = 4550.24
$ws.Range("I32").Value = 3705.4924
$ws.Range("J32").Value = 11625
$ws.Range("K32").Value = 3705.4924
$ws.Range("L32").Value = 11625
$ws.Range("M32").Value = -3418.4924
$ws.Range("N32").Value = -12199
$ws.Range("H45").Value = 1032.2
$ws.Range("I45").Value = 977.3570999999999
$ws.Range("J45").Value = 1800
$ws.Range("K45").Value = 977.3570999999999
$ws.Range("L45").Value = 1800
$ws.Range("M45").Value = -600.3570999999999
$ws.Range("N45").Value = -2554
$ws.Range("H74").Value = 694.43335
$ws.Range("I74").Value = 701.13794
$ws.Range("J74").Value = 500
$ws.Range("K74").Value = 701.13794
$ws.Range("L74").Value = 500
$ws.Range("M74").Value = 172.86206
$ws.Range("N74").Value = -2248
$ws.Range("H77").Value = 694.43335
$ws.Range("I77").Value = 701.13794
$ws.Range("J77").Value = 500
$ws.Range("K77").Value = 3505.6897
$ws.Range("L77").Value = 2500
$ws.Range("M77").Value = 862.3103000000001
$ws.Range("N77").Value = -11236
$ws.Range("H97").Value = 578.5
$ws.Range("I97").Value = 539.73334
$ws.Range("J97").Value = 1160
$ws.Range("K97").Value = 539.73334
$ws.Range("L97").Value = 1160
$ws.Range("M97").Value = -43.73334
$ws.Range("N97").Value = -2152
$ws.Range("H116").Value = 2833.3333
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 794
$ws.Range("N116").Value = -8088
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2833.3333
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 3500
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 3500
$ws.Range("M3").Value = -1386
$ws.Range("N3").Value = -3728
$ws.Range("H134").Value = 53134.45
$ws.Range("I134").Value = 86365.75
$ws.Range("J134").Value = 3287.5
$ws.Range("K134").Value = 259097.25
$ws.Range("L134").Value = 9862.5
$ws.Range("M134").Value = -256562.25
$ws.Range("N134").Value = -14932.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4275809
$ws.Range("I31").Value = 1470.125
$ws.Range("J31").Value = 23815644
$ws.Range("K31").Value = 1470.125
$ws.Range("L31").Value = 23815644
$ws.Range("M31").Value = -1175.125
$ws.Range("N31").Value = -23816234
$ws.Range("H34").Value = 4275809
$ws.Range("I34").Value = 1470.125
$ws.Range("J34").Value = 23815644
$ws.Range("K34").Value = 1470.125
$ws.Range("L34").Value = 23815644
$ws.Range("M34").Value = -1268.125
$ws.Range("N34").Value = -23816048
$ws.Range("H35").Value = 45114.223
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 50628.5
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 50628.5
$ws.Range("M35").Value = -706
$ws.Range("N35").Value = -51216.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 607.1429000000001
$ws.Range("J15").Value = 1500
$ws.Range("L15").Value = 4500
$ws.Range("N15").Value = -4780
$ws.Range("H68").Value = 332
$ws.Range("I68").Value = 354.8
$ws.Range("J68").Value = 275
$ws.Range("K68").Value = 1064.4
$ws.Range("L68").Value = 825
$ws.Range("M68").Value = -253.4000000000001
$ws.Range("N68").Value = -2447
$ws.Range("H71").Value = 332
$ws.Range("I71").Value = 354.8
$ws.Range("J71").Value = 275
$ws.Range("K71").Value = 3193.2
$ws.Range("L71").Value = 2475
$ws.Range("M71").Value = 862.7999999999997
$ws.Range("N71").Value = -10587
$ws.Range("H131").Value = 1854406.6
$ws.Range("I131").Value = 5001.25
$ws.Range("J131").Value = 3087343.8
$ws.Range("K131").Value = 15003.75
$ws.Range("L131").Value = 9262031.399999999
$ws.Range("M131").Value = -9963.75
$ws.Range("N131").Value = -9272111.399999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11337246
$ws.Range("I70").Value = 13080936
$ws.Range("J70").Value = 3266.6667
$ws.Range("K70").Value = 13080936
$ws.Range("L70").Value = 3266.6667
$ws.Range("M70").Value = -13080666
$ws.Range("N70").Value = -3806.6667
$ws.Range("H73").Value = 11337246
$ws.Range("I73").Value = 13080936
$ws.Range("J73").Value = 3266.6667
$ws.Range("K73").Value = 13080936
$ws.Range("L73").Value = 3266.6667
$ws.Range("M73").Value = -13080000
$ws.Range("N73").Value = -5138.6667
$ws.Range("H126").Value = 2224.2856
$ws.Range("J126").Value = 1224.5
$ws.Range("L126").Value = 3673.5
$ws.Range("N126").Value = -8613.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1001
$ws.Range("I32").Value = 1001
$ws.Range("K32").Value = 1001
$ws.Range("M32").Value = -684
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 2852.28
$ws.Range("J122").Value = 2875.625
$ws.Range("L122").Value = 8626.875
$ws.Range("N122").Value = -13526.875
$ws.Range("H132").Value = 6125.737
$ws.Range("I132").Value = 7741.9614
$ws.Range("J132").Value = 2623.9167
$ws.Range("K132").Value = 23225.8842
$ws.Range("L132").Value = 7871.750100000001
$ws.Range("M132").Value = -20695.8842
$ws.Range("N132").Value = -12931.7501
$ws.Range("H136").Value = 4183.4287
$ws.Range("I136").Value = 5083.304
$ws.Range("J136").Value = 2458.6667
$ws.Range("K136").Value = 15249.912
$ws.Range("L136").Value = 7376.000100000001
$ws.Range("M136").Value = -12699.912
$ws.Range("N136").Value = -12476.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 116395
$ws.Range("J98").Value = 116395
$ws.Range("L98").Value = 116395
$ws.Range("N98").Value = -122385
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 602
$ws.Range("J103").Value = 602
$ws.Range("L103").Value = 602
$ws.Range("N103").Value = -2946
$ws.Range("H136").Value = 4266
$ws.Range("I136").Value = 4810.129
$ws.Range("J136").Value = 2157.5
$ws.Range("K136").Value = 14430.387
$ws.Range("L136").Value = 6472.5
$ws.Range("M136").Value = -11880.387
$ws.Range("N136").Value = -11572.5
